{"js": "// Remove the paragraph:\n//   \"DI DARE ATTO che non sussistono oneri di sicurezza dovuti a\n//   rischio da interferenze;\"\n// which sits between the \"CAMPO.PROCEDERE\" bookmark paragraph and the\n// \"CAMPO.CONFERMARE\" bookmark paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"che non sussistono oneri di sicurezza dovuti a rischio da interferenze\";\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(needle) !== -1) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the paragraph:\n#   \"DI DARE ATTO che non sussistono oneri di sicurezza dovuti a\n#   rischio da interferenze;\"\n# which sits between the \"CAMPO.PROCEDERE\" bookmark paragraph and the\n# \"CAMPO.CONFERMARE\" bookmark paragraph.\n\n$d = $word.ActiveDocument\n\n$needle = \"che non sussistono oneri di sicurezza dovuti a rischio da interferenze\"\n\n# Collect the matching paragraph(s) first, then delete - mutating the\n# Paragraphs collection while iterating it is unsafe.\n$targets = New-Object System.Collections.ArrayList\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$needle*\") {\n        [void]$targets.Add($p.Range)\n    }\n}\n\nforeach ($r in $targets) {\n    $r.Delete()\n}\n"}
